# Smoke QA Mayo R33 - modificación de Data
$wb = $excel.ActiveWorkbook

# --- DatosCuenta sheet ---
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "SmokeMayo"
$wsCuenta.Range("B2").Value = "SmokeMayoLastName"
$wsCuenta.Range("C2").Value = 27100107
$wsCuenta.Range("D2").Value = 109
$wsCuenta.Activate()
$wsCuenta.Range("D14").Select()

# --- DatosHogar sheet ---
$wsHogar = $wb.Worksheets.Item("DatosHogar")
$wsHogar.Range("A2").Value = 629

# --- DatosMotor sheet ---
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Range("A2").Value = "SMA010"
$wsMotor.Range("B2").Value = "ABC12SSMA010"
$wsMotor.Range("C2").Value = "ZAZ123SSMA010"

# --- DatosAP sheet ---
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Activate()
$wsAP.Range("A2").Value = 21200109
$wsAP.Range("A3").Select()
